# Austrian population statistics workbook update
# - Add a "Number of districts per province" table (rows 9-26) to the
#   Sheet1 tab (sheet2.xml), counting districts per province via COUNTIFS
#   against the Foglio1 data.
# - Leaves a view/selection state similar to the one left behind after the
#   edit (scrolled down in Foglio1, selection left on the new table).

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Foglio1")
$ws  = $wb.Worksheets.Item("Sheet1")

# Header row for the new "Number of districts" table
$ws.Range("A9").Value = "Province"
$ws.Range("B9").Value = "Number of districts"

# One row per province, counting how many district rows in Foglio1
# (A2:A182) belong to it.
$provinces = @(
    "Milano",
    "Bergamo",
    "Brescia",
    "Como",
    "Cremona",
    "Lodi e Crema",
    "Mantova",
    "Pavia",
    "Sondrio",
    "Venezia ",
    "Belluno",
    "Padova   ",
    "Rovigo",
    "Treviso",
    "Udine  ",
    "Verona ",
    "Vicenza     "
)

$row = 10
foreach ($prov in $provinces) {
    $ws.Range("A$row").Value = $prov
    $ws.Range("B$row").Formula = '=COUNTIFS(Foglio1!$A$2:$A$182,Sheet1!$A' + $row + ')'
    $row = $row + 1
}

# Recalculate so the cached <v> values are correct in the saved file.
$excel.Calculate()

# --- View state -----------------------------------------------------
# Scroll Foglio1 down and select the whole Province column, as if the
# user had just been reading off values to build the table above.
[void]$src.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$src.Range("A2").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 148
$src.Range("A1:A182").Select() | Out-Null

# Leave the new table's sheet active with the selection where editing
# stopped.
[void]$ws.Activate()
$ws.Range("F16").Select() | Out-Null
